# Commit message summary: parser/importer refactor added tests that read
# individual sheets by sheetNumber, so the mock workbook grew two extra
# sheets ("Sheet3" then "Sheet2", in that tab order after "Sheet1") that
# carry the same table, just reshuffled, plus "Sheet3" became the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create "Sheet3" right after "Sheet1" ------------------------------
# Copy()ing Sheet1 clones its cell data, styles, row heights and page
# setup, so the clone starts out identical and we only touch the cells
# that actually differ in the target sheet.
$ws1.Copy($null, $ws1)
$ws3 = $wb.Worksheets.Item(2)
$ws3.Name = "Sheet3"

# --- Create "Sheet2" right after "Sheet3" -------------------------------
$ws3.Copy($null, $ws3)
$ws2 = $wb.Worksheets.Item(3)
$ws2.Name = "Sheet2"

# --- Reorder the data rows on "Sheet3" ----------------------------------
$ws3.Range("A2").Value = "São Paulo"
$ws3.Range("B2").Value = "SP"
$ws3.Range("C2").Value = 41644
$ws3.Range("E2").Value = 14

$ws3.Range("A3").Value = "Campo Grande"
$ws3.Range("B3").Value = "MS"
$ws3.Range("C3").Value = 41641
$ws3.Range("E3").Value = 11

$ws3.Range("A4").Value = "Rio de Janeiro"
$ws3.Range("B4").Value = "RJ"
$ws3.Range("C4").Value = 41642
$ws3.Range("E4").Value = 12

$ws3.Range("A5").Value = "São Paulo"
$ws3.Range("B5").Value = "SP"
$ws3.Range("C5").Value = 41643
$ws3.Range("E5").Value = 13

$ws3.Range("A6").Value = "Pernambuco"
$ws3.Range("B6").Value = "PE"
$ws3.Range("C6").Value = 41640
$ws3.Range("E6").Value = 10

# "Sheet3" also uses a slightly taller default row height (12.8 vs 12.1)
$ws3.UsedRange.RowHeight = 12.8

# --- Reorder the data rows on "Sheet2" (descending by Número) ----------
$ws2.Range("A2").Value = "São Paulo"
$ws2.Range("B2").Value = "SP"
$ws2.Range("C2").Value = 41644
$ws2.Range("E2").Value = 14

$ws2.Range("A3").Value = "São Paulo"
$ws2.Range("B3").Value = "SP"
$ws2.Range("C3").Value = 41643
$ws2.Range("E3").Value = 13

$ws2.Range("A4").Value = "Rio de Janeiro"
$ws2.Range("B4").Value = "RJ"
$ws2.Range("C4").Value = 41642
$ws2.Range("E4").Value = 12

$ws2.Range("A5").Value = "Campo Grande"
$ws2.Range("B5").Value = "MS"
$ws2.Range("C5").Value = 41641
$ws2.Range("E5").Value = 11

$ws2.Range("A6").Value = "Pernambuco"
$ws2.Range("B6").Value = "PE"
$ws2.Range("C6").Value = 41640
$ws2.Range("E6").Value = 10

# --- Per-sheet selection, matching each sheet's saved cursor -----------
$ws1.Range("A1").Select()
$ws2.Range("A5").Select()
$ws3.Range("A7").Select()

# --- Active tab / window state ------------------------------------------
# "Sheet3" (the 2nd tab) becomes the selected/active sheet.
$ws3.Activate()
$ws3.Range("A7").Select()

# tabRatio (the tab-bar/horizontal-scrollbar split) - best effort; not all
# hosts persist this window chrome setting back into workbookView.
try { $excel.ActiveWindow.TabRatio = 0.567 } catch {}
